# Swap the deck's colour theme: the presentation currently uses the
# "Integral" / Red Violet colour scheme (the theme linked from the slide
# master). The target state is the default "Office" colour scheme - the
# font scheme and the fill/line/effect format scheme are already
# identical between the two themes in this deck, so re-pointing every
# one of the twelve theme colour slots is all that's needed to turn the
# active theme into the Office colour scheme.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      -> 000000
$cs.Colors(2).RGB  = 16777215   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      -> 44546A
$cs.Colors(4).RGB  = 15132391   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  -> ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  -> FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  -> 4472C4
$cs.Colors(10).RGB = 4697456    # accent6  -> 70AD47
$cs.Colors(11).RGB = 12673797   # hlink    -> 0563C1
$cs.Colors(12).RGB = 7491477    # folHlink -> 954F72
